$d = $word.ActiveDocument

# Locate the sentence containing the entry-fee deadline date so we do not
# depend on hard-coded character offsets.
$dateRng = $d.Content
$found = $dateRng.Find.Execute("Friday 25th February.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Friday 25th February.' sentence to edit"
}
$sentenceStart = $dateRng.Start
$sentenceEnd = $dateRng.End

# Positions (relative to $sentenceStart) of the pieces of the sentence:
#   "Friday " (7 chars) "25" (2 chars) "th" (2 chars) " February" (9 chars) "." (1 char)
$dayStart = $sentenceStart + 7          # start of "25"
$splitPos = $dayStart + 1               # between "2" and "5"
$dotStart = $sentenceEnd - 1            # start of the trailing "."
$dotEnd = $sentenceEnd                  # end of the trailing "."

# The trailing "." run shares identical run formatting with the preceding
# " February" run. Editing text elsewhere in this paragraph causes the
# runtime to re-merge all adjacent, identically-formatted runs, which would
# incorrectly fuse " February" and "." together. Give the "." run a
# momentarily distinct font size so it keeps its own run, then restore the
# original size once the real edit has been made.
$dotRng = $d.Range($dotStart, $dotEnd)
$originalDotSize = $dotRng.Font.Size
$dotRng.Font.Size = $originalDotSize + 1

# Insert the "_GoBack" bookmark between the "2" and "5" of "25" -- this is
# the point where the document was last edited, which forces a run split at
# that location (matching how Word splits runs around an edit point).
$bmRng = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# Change "5" -> "4" so the date becomes "24th February."
$digitRng = $d.Range($splitPos, $splitPos + 1)
$digitRng.Text = "4"

# Restore the trailing "." run back to its original formatting.
$dotRng2 = $d.Range($dotStart, $dotEnd)
$dotRng2.Font.Size = $originalDotSize
